# "Finalizo las pruebas del ejercicio del algoritmo que coordina."
#
# - Renames "Creacion de un nuevo proyecto" -> "Genere el proyecto por crear"
#   (kept in place so its comments/table/drawing stay attached).
# - Inserts a brand new sheet "Valide el nuevo proyecto" right before it,
#   with the Gherkin scenario + examples tables that validate a new project.
# - Restores view-state (selections / active tab) to match the target.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing sheet in place (preserves its comments/vmlDrawing/table)
# ---------------------------------------------------------------------------
$genere = $wb.Worksheets.Item("Creacion de un nuevo proyecto")
$genere.Name = "Genere el proyecto por crear"

# Donor cells (already-existing styles) we copy *formats only* from, so the
# new cells land on the very same style index Excel itself would reuse.
$donorDate = $genere.Range("B9")          # numFmtId 14 "m/d/yyyy"

# ---------------------------------------------------------------------------
# 2. Insert the new sheet right before it
# ---------------------------------------------------------------------------
$valide = $wb.Worksheets.Add($genere, $null)
$valide.Name = "Valide el nuevo proyecto"

# ---------------------------------------------------------------------------
# 3. Content
# ---------------------------------------------------------------------------
$valide.Range("A1").Value = "Escenario"
$valide.Range("A1").Style = "Accent2"
$valide.Range("B1").Value = "Valide el nuevo proyecto"

$valide.Range("A3").Value = "Dados"
$valide.Range("A3").Font.Bold = $true
$valide.Range("B3").Value = "los proyectos existentes"

$valide.Range("B4").Value = "Nombre"
$valide.Range("B5").Value = "Proyecto 1"
$valide.Range("B6").Value = "Proyecto 2"

$valide.Range("A8").Value = "Cuando"
$valide.Range("A8").Font.Bold = $true
$valide.Range("B8").Value = "se valida el nuevo proyecto con nombre <nombre> y la fecha de inicio es <fecha de inicio>"

$valide.Range("A9").Value = "Entonces"
$valide.Range("A9").Font.Bold = $true
$valide.Range("B9").Value = "se indica si <es valido o no>"

$valide.Range("A11").Value = "Ejemplos"
$valide.Range("A11").Font.Bold = $true

$valide.Range("A12").Value = "Ejemplo"
$valide.Range("A12").Font.Bold = $true
$valide.Range("B12").Value = "Nombre"
$valide.Range("C12").Value = "Fecha de inicio"
$valide.Range("D12").Value = "Es valido o no"

$valide.Range("A13").Value = "Son validos"
$valide.Range("B13").Value = "Proyecto asombroso"
$valide.Range("C13").Value = 43399
$donorDate.Copy()
$valide.Range("C13").PasteSpecial(-4122)
$valide.Range("D13").Value = "Es válido"

$valide.Range("A14").Value = "El nombre es requerido"
$valide.Range("C14").Value = 43399
$donorDate.Copy()
$valide.Range("C14").PasteSpecial(-4122)
$valide.Range("D14").Value = "Es inválido"

$valide.Range("A15").Value = "El nombre tiene un tamaño máximo de 200 caracteres"
$valide.Range("B15").Value = '"1234567890123456789012345678901234567890123456789012345678901234567890123456789012345678901234567890123456789012345678901234567890123456789012345678901234567890123456789012345678901234567890123456789012345678901"'
$valide.Range("C15").Value = 43399
$donorDate.Copy()
$valide.Range("C15").PasteSpecial(-4122)
$valide.Range("D15").Value = "Es inválido"

$valide.Range("A16").Value = "El nombre no puede contener solamente espacios en blanco"
$valide.Range("B16").Value = '"    "'
$valide.Range("C16").Value = 43399
$donorDate.Copy()
$valide.Range("C16").PasteSpecial(-4122)
$valide.Range("D16").Value = "Es inválido"

$valide.Range("A17").Value = "El nombre sólo puede contener letras y números"
$valide.Range("B17").Value = "Hola-Mundo"
$valide.Range("C17").Value = 43399
$donorDate.Copy()
$valide.Range("C17").PasteSpecial(-4122)
$valide.Range("D17").Value = "Es inválido"

$valide.Range("A18").Value = "El nombre debe ser único"
$valide.Range("B18").Value = "Proyecto 1"
$valide.Range("C18").Value = 43399
$donorDate.Copy()
$valide.Range("C18").PasteSpecial(-4122)
$valide.Range("D18").Value = "Es inválido"

$valide.Range("A19").Value = "La fecha de inicio es requerida"
$valide.Range("B19").Value = "Proyecto asombroso"
$valide.Range("C19").Value = "nula"
$valide.Range("D19").Value = "Es inválido"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Column widths (best effort - engine quantises to 1/6 character steps)
# ---------------------------------------------------------------------------
$valide.Columns.Item(1).ColumnWidth = 54.83
$valide.Columns.Item(2).ColumnWidth = 27.67
$valide.Columns.Item(3).ColumnWidth = 15.67
$valide.Columns.Item(4).ColumnWidth = 14.67

# ---------------------------------------------------------------------------
# 5. Tables
# ---------------------------------------------------------------------------
$tblEjemplos = $valide.ListObjects.Add(1, $valide.Range("A12:D19"), $null, 1)
$tblEjemplos.Name = "Table2"
$tblEjemplos.TableStyle = "TableStyleLight9"

$tblNombre = $valide.ListObjects.Add(1, $valide.Range("B4:B6"), $null, 1)
$tblNombre.Name = "Table3"
$tblNombre.TableStyle = "TableStyleLight9"

$valide.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 6. View state: selections + active tab
# ---------------------------------------------------------------------------
$valide.Range("A13:A19").Select()

$genere.Range("A1:B10").Select()
$genere.Activate()
